# Swap the presentation's theme colour scheme from "Integral" to the
# stock "Office Theme" palette (the slide master's theme -> ppt/theme/theme1.xml).
# The font scheme (Arial majors/minors) and format scheme (fills/lines/
# effects) are identical between the two themes, so only the 12 theme
# colour slots need to change.
#
# COM colour-slot order (matches PowerPoint's ColorScheme.Colors index):
#   1 dk1, 2 lt1, 3 dk2, 4 lt2,
#   5 accent1, 6 accent2, 7 accent3, 8 accent4, 9 accent5, 10 accent6,
#   11 hlink, 12 folHlink
# RGB is passed as a COM BGR-packed long (0xBBGGRR), matching VBA's RGB().

$p = $ppt.ActivePresentation
$cs = $p.SlideMaster.Theme.ThemeColorScheme

function BGR($r, $g, $b) {
    return ($b * 65536) + ($g * 256) + $r
}

$cs.Colors(1).RGB  = BGR 0x00 0x00 0x00   # dk1      000000
$cs.Colors(2).RGB  = BGR 0xFF 0xFF 0xFF   # lt1      FFFFFF
$cs.Colors(3).RGB  = BGR 0x44 0x54 0x6A   # dk2      44546A
$cs.Colors(4).RGB  = BGR 0xE7 0xE6 0xE6   # lt2      E7E6E6
$cs.Colors(5).RGB  = BGR 0x5B 0x9B 0xD5   # accent1  5B9BD5
$cs.Colors(6).RGB  = BGR 0xED 0x7D 0x31   # accent2  ED7D31
$cs.Colors(7).RGB  = BGR 0xA5 0xA5 0xA5   # accent3  A5A5A5
$cs.Colors(8).RGB  = BGR 0xFF 0xC0 0x00   # accent4  FFC000
$cs.Colors(9).RGB  = BGR 0x44 0x72 0xC4   # accent5  4472C4
$cs.Colors(10).RGB = BGR 0x70 0xAD 0x47   # accent6  70AD47
$cs.Colors(11).RGB = BGR 0x05 0x63 0xC1   # hlink    0563C1
$cs.Colors(12).RGB = BGR 0x95 0x4F 0x72   # folHlink 954F72
